# Applies the "excluded block of text is created" edit:
#  1. Move the _GoBack bookmark from the end of the "Martindale..." reference
#     paragraph to the blank paragraph right before "COMPARISON OF PLE AND LMS".
#  2. Merge the two adjacent runs ("...for themselves" + ". (7 ") into a
#     single run so the sentence is no longer split mid-way.
#  3. Turn the list-numbered paragraph that only holds the page break into a
#     plain paragraph, and add a new blank plain paragraph right before it
#     (the "excluded block of text" -- an empty paragraph -- gets created).

$d = $word.ActiveDocument

# --- Step 1: merge the split sentence back into a single run -------------
$findRange = $d.Content
$null = $findRange.Find.Execute("for themselves")
$boundary = $findRange.End
$tail = $d.Range($boundary, $boundary + 5)
if ($tail.Text -eq ". (7 ") {
    $tail.Delete()
    $insertionPoint = $d.Range($boundary, $boundary)
    $insertionPoint.InsertAfter(". (7 ")
}

# --- Step 2: relocate the _GoBack bookmark --------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$refRange = $d.Content
$null = $refRange.Find.Execute("COMPARISON OF PLE AND LMS")
$precedingPara = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -eq $refRange.Start) {
        $precedingPara = $d.Paragraphs.Item($i - 1)
        break
    }
}
$d.Bookmarks.Add("_GoBack", $precedingPara.Range) | Out-Null

# --- Step 3: clear the list formatting off the page-break paragraph, and --
# --- insert a new blank paragraph in front of it --------------------------
# Locate the (list-numbered) paragraph whose only content is a page break --
# it immediately follows the "Martindale..., 177-193." reference paragraph.
$pageBreakPara = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.ListFormat.ListType -ne 0 -and $pp.Range.Text -eq ([string][char]12 + [string][char]13)) {
        $pageBreakPara = $pp
        break
    }
}
$pageBreakStart = $pageBreakPara.Range.Start

# Replace the paragraph's XML to drop the List Paragraph style / numbering,
# keeping its run (the page break) and the paragraph-mark language setting.
$cleanXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
$pageBreakPara.Range.Duplicate().InsertXML($cleanXml)

# Insert a brand-new, empty plain paragraph right before it.
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$d.Range($pageBreakStart, $pageBreakStart).InsertXML($blankXml)

Write-Output "edit complete"
